$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 165 (achievement id 263) ---
# Write string-valued cells first, in the exact order new shared strings
# need to be allocated: E, C, B.
$ws.Range("E165").Value = "TLT?[2036]"
$ws.Range("C165").Value = "自 投 罗 网"
$ws.Range("B165").Value = "经 典 老 歌"
$ws.Range("A165").Value = 263
$ws.Range("D165").Value = 2
$ws.Range("F165").Value = 1
$ws.Range("G165").Value = "START"

# --- Row 166 (achievement id 264) ---
$ws.Range("C166").Value = "集齐四大悲剧"
$ws.Range("B166").Value = "莎比"
$ws.Range("E166").Value = "(ATLT?[2028])&(ATLT?[2029])&(ATLT?[2030])&(ATLT?[2031])"
$ws.Range("A166").Value = 264
$ws.Range("D166").Value = 2
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = "START"

# --- Row 167 (achievement id 265) ---
# Column E has no sheet-level default style, so a brand-new E167 cell
# would otherwise come out unstyled; borrow E166's formatting first so it
# matches the rest of the table (s="5") without minting a new style.
$ws.Range("E166").Copy()
$ws.Range("E167").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B167").Value = "死了但没完全死"
$ws.Range("C167").Value = "死而复生"
$ws.Range("E167").Value = "EVT?[20000,20001,11504]"
$ws.Range("A167").Value = 265
$ws.Range("D167").Value = 2
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = "TRAJECTORY"

# --- Row 168: trailing blank spacer row (matches prior blank rows 165/166) ---
# Carry over the same "styled but empty" E-column cell the sheet always
# keeps as a spacer below the last data row, by copying formatting only.
$ws.Range("E164").Copy()
$ws.Range("E168").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update achievement 219 (row 121) condition to add the new ATLT clause ---
$ws.Range("E121").Value = "(ATLT?[1023])&(ATLT?[1048])&(ATLT?[1064])&(ATLT?[1114])&(ATLT?[1135])&(ATLT?[1141])&(ATLT?[1147])"

# --- Leave the selection where the editor ended up working ---
$ws.Activate()
$ws.Range("G122").Select()
